$d = $word.ActiveDocument
$find = $d.Content.Find
$find.ClearFormatting()
$find.Replacement.ClearFormatting()
$find.Execute("Are there any syntax errors in the code?", $true, $false, $false, $false, $false, $true, 1, $false, "Is the code free of syntax errors?", 2)
